# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F value for the "展览" sheet
$exhibitionUpdates = @{
    2  = 155
    3  = 1718
    4  = 791
    5  = 1124
    7  = 12001
    10 = 478
    13 = 865
    14 = 13483
    15 = 13479
    16 = 40
    20 = 486
    23 = 966
    24 = 175
}

# Row -> new F value for the "全部类型" sheet (values differ slightly from
# "展览" for row 23, same as observed in the source data)
$allTypesUpdates = @{
    2  = 155
    3  = 1718
    4  = 791
    5  = 1124
    7  = 12001
    10 = 478
    13 = 865
    14 = 13483
    15 = 13479
    16 = 40
    20 = 486
    23 = 968
    24 = 175
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
